$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.268731832504272
$ws.Range("B1").Value = 2.719936370849609
$ws.Range("C1").Value = 4.910079479217529
$ws.Range("D1").Value = 2.034052848815918
$ws.Range("E1").Value = 1.034162402153015
